$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$ws.Range("A2").Value = "c36c75c"
$ws.Range("N2").Value = "CATS"
$ws.Range("A3").Value = "092adc4"
$ws.Range("M3").Value = "japanese"
$ws.Range("A4").Value = "210be5a"
$ws.Range("M4").Value = "japanese"
$ws.Range("N4").Value = "CATS"
$ws.Range("A5").Value = "1e90dc1"
$ws.Range("M5").Value = "japanese"
$ws.Range("N5").Value = "CATS"
$ws.Range("A6").Value = "b71efee"
$ws.Range("N6").Value = "DOGS"
$ws.Range("A7").Value = "2109cf6"
$ws.Range("N7").Value = "FISH"
$ws.Range("A8").Value = "8ffc71d"
$ws.Range("N8").Value = "BIRDS"
$ws.Range("A9").Value = "1aa5347"
$ws.Range("N9").Value = "CATS"
$ws.Range("A10").Value = "e794a78"
$ws.Range("N10").Value = "DOGS"
$ws.Range("A11").Value = "4d49497"
$ws.Range("N11").Value = "BIRDS"
$ws.Range("A12").Value = "d9cf7b8"
$ws.Range("N12").Value = "REPTILES"
$ws.Range("A13").Value = "0da479a"
$ws.Range("A14").Value = "'1563471"
$ws.Range("M14").Value = "japanese"
$ws.Range("N14").Value = "DOGS"
$ws.Range("A15").Value = "cea858b"
$ws.Range("N15").Value = "BIRDS"
$ws.Range("A16").Value = "2ebb06c"
$ws.Range("N16").Value = "FISH"
$ws.Range("A17").Value = "53316bd"
$ws.Range("N17").Value = "DOGS"
$ws.Range("A18").Value = "e190b2e"
$ws.Range("N18").Value = "REPTILES"
$ws.Range("A19").Value = "'1449692"
$ws.Range("M19").Value = "japanese"
$ws.Range("A20").Value = "2b2f079"
$ws.Range("N20").Value = "BIRDS"
$ws.Range("A21").Value = "b6ca709"
$ws.Range("N21").Value = "REPTILES"
$ws.Range("A22").Value = "d28a183"
$ws.Range("N22").Value = "DOGS"
$ws.Range("A23").Value = "bbde1e0"
$ws.Range("M23").Value = "japanese"
$ws.Range("N23").Value = "DOGS"
$ws.Range("A24").Value = "ff6247c"
$ws.Range("M24").Value = "english"
$ws.Range("N24").Value = "DOGS"
$ws.Range("A25").Value = "3f9ea5e"
$ws.Range("N25").Value = "DOGS"
$ws.Range("A26").Value = "ef6dac0"
$ws.Range("N26").Value = "BIRDS"
$ws.Range("A27").Value = "e5d3fce"
$ws.Range("M27").Value = "english"
$ws.Range("N27").Value = "FISH"
$ws.Range("A28").Value = "a046401"
$ws.Range("M28").Value = "japanese"
$ws.Range("N28").Value = "CATS"
$ws.Range("A29").Value = "eb0c5b1"
$ws.Range("M29").Value = "japanese"
$ws.Range("A30").Value = "66e82a2"
$ws.Range("N30").Value = "DOGS"
$ws.Range("A31").Value = "c86a4e2"
$ws.Range("M31").Value = "japanese"
$ws.Range("N31").Value = "BIRDS"
$ws.Range("A32").Value = "abd949b"
$ws.Range("M32").Value = "english"
$ws.Range("A33").Value = "f7f0091"
$ws.Range("M33").Value = "japanese"
$ws.Range("N33").Value = "BIRDS"
$ws.Range("A34").Value = "c4fc0c8"
$ws.Range("A35").Value = "78e02fc"
$ws.Range("M35").Value = "japanese"
$ws.Range("A36").Value = "1d3dd2f"
$ws.Range("M36").Value = "japanese"
$ws.Range("N36").Value = "DOGS"
$ws.Range("A37").Value = "dce69c2"
$ws.Range("N37").Value = "DOGS"
$ws.Range("A38").Value = "6cb08c9"
$ws.Range("M38").Value = "english"
$ws.Range("N38").Value = "BIRDS"
$ws.Range("A39").Value = "1e4d2c6"
$ws.Range("M39").Value = "japanese"
$ws.Range("N39").Value = "BIRDS"
$ws.Range("A40").Value = "9385d98"
$ws.Range("M40").Value = "english"
$ws.Range("N40").Value = "REPTILES"
$ws.Range("A41").Value = "5082b96"
$ws.Range("M41").Value = "english"
$ws.Range("N41").Value = "REPTILES"
$ws.Range("A42").Value = "'7585309"
$ws.Range("M42").Value = "english"
$ws.Range("N42").Value = "BIRDS"
$ws.Range("A43").Value = "'5437757"
$ws.Range("M43").Value = "japanese"
$ws.Range("A44").Value = "7722a75"
$ws.Range("M44").Value = "english"
$ws.Range("N44").Value = "REPTILES"
$ws.Range("A45").Value = "0cf1d93"
$ws.Range("N45").Value = "DOGS"
$ws.Range("A46").Value = "a0bbc2c"
$ws.Range("M46").Value = "english"
$ws.Range("N46").Value = "FISH"
$ws.Range("A47").Value = "f1f15dd"
$ws.Range("N47").Value = "CATS"
$ws.Range("A48").Value = "d5dce40"
$ws.Range("M48").Value = "english"
$ws.Range("N48").Value = "CATS"
$ws.Range("A49").Value = "b5601e4"
$ws.Range("A50").Value = "a123e52"
$ws.Range("N50").Value = "REPTILES"
$ws.Range("A51").Value = "a49811b"
$ws.Range("M51").Value = "japanese"
$ws.Range("N51").Value = "DOGS"

Write-Output "Applied users sheet updates."
